$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3157
$ws.Range("J3").Value = 3297
$ws.Range("I4").Value = 1761
$ws.Range("J4").Value = 739
$ws.Range("J5").Value = 257
$ws.Range("J6").Value = 3899
$ws.Range("I7").Value = 26207
$ws.Range("J7").Value = 11349

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 35
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 363

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 55
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 169
$ws.Range("J7").Value = 404

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 31
$ws.Range("J7").Value = 343
$ws.Range("J8").Value = 725
$ws.Range("J9").Value = 68
$ws.Range("J10").Value = 67
$ws.Range("J11").Value = 163
$ws.Range("J12").Value = 22
$ws.Range("J14").Value = 47
$ws.Range("J15").Value = 132
$ws.Range("J20").Value = 239
$ws.Range("J25").Value = 63
$ws.Range("J27").Value = 68
$ws.Range("J29").Value = 651
$ws.Range("J30").Value = 47
$ws.Range("J33").Value = 490
$ws.Range("J34").Value = 57
$ws.Range("J37").Value = 363
$ws.Range("J42").Value = 453
$ws.Range("J44").Value = 88
$ws.Range("J46").Value = 40
$ws.Range("J47").Value = 85
$ws.Range("J51").Value = 151
$ws.Range("J52").Value = 307
$ws.Range("J53").Value = 110
$ws.Range("J54").Value = 215
$ws.Range("J55").Value = 147
$ws.Range("J57").Value = 51
$ws.Range("J60").Value = 78
$ws.Range("I63").Value = 216
$ws.Range("J63").Value = 51
$ws.Range("J64").Value = 77
$ws.Range("J67").Value = 404
$ws.Range("J71").Value = 41
$ws.Range("J72").Value = 45
$ws.Range("J75").Value = 35
$ws.Range("J77").Value = 97
$ws.Range("J78").Value = 152
$ws.Range("J79").Value = 338
$ws.Range("J83").Value = 259
$ws.Range("J85").Value = 518
$ws.Range("J88").Value = 117
$ws.Range("J89").Value = 130
$ws.Range("J91").Value = 127
$ws.Range("J92").Value = 35
$ws.Range("J95").Value = 177
$ws.Range("J99").Value = 165
$ws.Range("I101").Value = 26207
$ws.Range("J101").Value = 11349

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 99
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 177

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 130
$ws.Range("J3").Value = 156
$ws.Range("J6").Value = 161
$ws.Range("J7").Value = 490

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 196
$ws.Range("J3").Value = 223
$ws.Range("J6").Value = 166
$ws.Range("J7").Value = 651

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 127
$ws.Range("J3").Value = 195
$ws.Range("J4").Value = 40
$ws.Range("J7").Value = 518

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 93
$ws.Range("J3").Value = 99
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 228
$ws.Range("J7").Value = 453

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 17
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 52
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 37
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 122
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 74
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 69
$ws.Range("J3").Value = 87
$ws.Range("J6").Value = 134
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 36
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 218
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 217
$ws.Range("J7").Value = 725

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 39
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 114
$ws.Range("J3").Value = 103
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 343

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 22
